$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2,3,5,6,7 are cyclically permuted (all columns except D,L,M,N,O,P,R,S
# stay identical across these rows, so only those columns need updating).
#
# Target values (post-edit), derived from the diff:
#   Row 2: D=44252 L=Primera M=120 N=13000 O=14000 P=13500 R=Región Metropolitana      S=750
#   Row 3: D=44250 L=Primera M=200 N=14000 O=15000 P=14500 R=Región Metropolitana      S=806
#   Row 5: D=44257 L=Primera M=100 N=14000 O=15000 P=14500 R=Región Metropolitana      S=806
#   Row 6: D=45072 L=Segunda M=100 N=16000 O=16000 P=16000 R=Provincia de Chacabuco    S=889
#   Row 7: D=45072 L=Segunda M=100 N=17000 O=17000 P=17000 R=Provincia de Limarí       S=944

$updates = @{
    2 = @{ D = 44252; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; R = "Región Metropolitana"; S = 750 }
    3 = @{ D = 44250; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana"; S = 806 }
    5 = @{ D = 44257; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana"; S = 806 }
    6 = @{ D = 45072; L = "Segunda"; M = 100; N = 16000; O = 16000; P = 16000; R = "Provincia de Chacabuco"; S = 889 }
    7 = @{ D = 45072; L = "Segunda"; M = 100; N = 17000; O = 17000; P = 17000; R = "Provincia de Limarí"; S = 944 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
